# Apply the property-land (issue #5) edits to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "土地" (land) - rename header, clean up text, add new
# property/legislator metadata columns I..O.
# ---------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

# Header text column A label changes from "土地坐落" to the raw field name.
$land.Cells.Item(1, 2).Value = "name"

# Clean up stray dashes/spaces in the land-parcel descriptions.
$land.Cells.Item(2, 2).Value = "桃園縣中壢市石頭段00490032地號"
$land.Cells.Item(3, 2).Value = "桃園縣中壢市石頭段00490033地號"
$land.Cells.Item(4, 2).Value = "桃園縣中壢市石頭段01320055地號"
$land.Cells.Item(5, 2).Value = "臺北市大安區大安段二小段01110000地號"

# Clean up stray space in the share-portion fraction.
$land.Cells.Item(5, 4).Value = "10000分之566"

# Copy the formats of the last existing column (H) onto the new columns
# (I..O) for both the header row and the data rows.
$land.Range("H1").Copy() | Out-Null
$land.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$land.Range("H2:H5").Copy() | Out-Null
$land.Range("I2:O5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New header row (I1:O1).
$land.Cells.Item(1, 9).Value  = "property_category"
$land.Cells.Item(1, 10).Value = "category"
$land.Cells.Item(1, 11).Value = "date"
$land.Cells.Item(1, 12).Value = "legislator_name"
$land.Cells.Item(1, 13).Value = "legislator_id"
$land.Cells.Item(1, 14).Value = "source_file"
$land.Cells.Item(1, 15).Value = "index"

# Row 2 (index 13).
$land.Cells.Item(2, 9).Value  = "land"
$land.Cells.Item(2, 10).Value = "normal"
$land.Cells.Item(2, 11).Value = "2012-04-20"
$land.Cells.Item(2, 12).Value = "廖正井"
$land.Cells.Item(2, 13).Value = 1711
$land.Cells.Item(2, 14).Value = "tmp845a1"
$land.Cells.Item(2, 15).Value = 13

# Row 3 (index 14).
$land.Cells.Item(3, 9).Value  = "land"
$land.Cells.Item(3, 10).Value = "normal"
$land.Cells.Item(3, 11).Value = "2012-04-20"
$land.Cells.Item(3, 12).Value = "廖正井"
$land.Cells.Item(3, 13).Value = 1711
$land.Cells.Item(3, 14).Value = "tmp845a1"
$land.Cells.Item(3, 15).Value = 14

# Row 4 (index 15).
$land.Cells.Item(4, 9).Value  = "land"
$land.Cells.Item(4, 10).Value = "normal"
$land.Cells.Item(4, 11).Value = "2012-04-20"
$land.Cells.Item(4, 12).Value = "廖正井"
$land.Cells.Item(4, 13).Value = 1711
$land.Cells.Item(4, 14).Value = "tmp845a1"
$land.Cells.Item(4, 15).Value = 15

# Row 5 (index 16).
$land.Cells.Item(5, 9).Value  = "land"
$land.Cells.Item(5, 10).Value = "normal"
$land.Cells.Item(5, 11).Value = "2012-04-20"
$land.Cells.Item(5, 12).Value = "廖正井"
$land.Cells.Item(5, 13).Value = 1711
$land.Cells.Item(5, 14).Value = "tmp845a1"
$land.Cells.Item(5, 15).Value = 16

# ---------------------------------------------------------------
# Sheet "建物" (building) - clean up stray dashes/spaces in the
# building descriptions / parking-space note.
# ---------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Cells.Item(2, 2).Value = "桃園縣中壢市石頭段06151000建號"
$building.Cells.Item(3, 2).Value = "桃圜縣中壢市石頭段06152000建號"
$building.Cells.Item(4, 2).Value = "臺北市大安區大安段二小段01874000建號"
$building.Cells.Item(5, 2).Value = "臺北市大安區大安段二小段01875000建號"
$building.Cells.Item(5, 8).Value = "(超過五年停車位）"

# ---------------------------------------------------------------
# Sheet "存款" (deposits) - remove stray internal spaces from bank
# branch names.
# ---------------------------------------------------------------
$deposit = $wb.Worksheets.Item("存款")
$deposit.Cells.Item(6, 2).Value  = "台北富邦商業銀行市府分行"
$deposit.Cells.Item(7, 2).Value  = "台北富邦商業銀行城中分行"
$deposit.Cells.Item(8, 2).Value  = "中華郵政股份有限公司桃圜府前郵局"
$deposit.Cells.Item(9, 2).Value  = "中華郵政股份有限公司台北信維郵局"
$deposit.Cells.Item(10, 2).Value = "中華郵政股份有限公司台北信維郵局"
$deposit.Cells.Item(11, 2).Value = "國泰世華商業銀行信義分行"
$deposit.Cells.Item(12, 2).Value = "國泰世華商業銀行信義分行"
$deposit.Cells.Item(13, 2).Value = "國泰世華商業銀行信義分行"
$deposit.Cells.Item(17, 2).Value = "中國信託商業銀行敦南分行"
$deposit.Cells.Item(18, 2).Value = "台北富邦商業銀行敦和分行"

# ---------------------------------------------------------------
# Sheet "事業投資" (business investment) - remove stray internal
# spaces from the investment address and acquisition date.
# ---------------------------------------------------------------
$invest = $wb.Worksheets.Item("事業投資")
$invest.Cells.Item(2, 4).Value = "臺北市中山區長春路378號6F"
$invest.Cells.Item(2, 6).Value = "100年07月14H"
